$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "contrast2"/"contrast3"/"contrast4" columns ---
# The old layout had a single "contrast" column at D; the new layout needs
# four (contrast1..contrast4) at D:G, so insert three blank columns after D.
$ws.Columns("E:G").Insert()

# --- Header row ---
$ws.Range("D1").Value = "contrast1"
$ws.Range("E1").Value = "contrast2"
$ws.Range("F1").Value = "contrast3"
$ws.Range("G1").Value = "contrast4"

# --- Fill in the blank contrast2/contrast3/contrast4 cells for existing rows ---
$ws.Range("E5").Value = "diff"
$ws.Range("F5").Value = "diff"
$ws.Range("G5").Value = "diff"

# --- New trial rows (6 through 13), continuing after the existing trial 5 ("Alert") ---

# Trial 6 - Sitagliptin
$ws.Range("B7").Value = "Sitagliptin"
$ws.Range("D7").Value = "diff"
$ws.Range("H7").Value = -0.001
$ws.Range("I7").Value = 0.012
$ws.Range("J7").Value = 0.013
$ws.Range("K7").Value = -0.004
$ws.Range("L7").Value = 0.013

# Trial 7 - Bamlanivimab/Etesevimab vs Sotrovimab
$ws.Range("B8").Value = "Bamlanivimab/Etesevimab"
$ws.Range("C8").Value = "Sotrovimab "
$ws.Range("D8").Value = "diff"
$ws.Range("H8").Value = 0.02
$ws.Range("I8").Value = 0.013
$ws.Range("J8").Value = 0.011
$ws.Range("K8").Value = 0.019
$ws.Range("L8").Value = 0.011
$ws.Range("M8").Value = -0.001
$ws.Range("N8").Value = 0.003
$ws.Range("O8").Value = 0.011
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0.011

# Trial 8 - Ivermectin
$ws.Range("B9").Value = "Ivermectin "
$ws.Range("D9").Value = "diff"
$ws.Range("H9").Value = 0.068
$ws.Range("I9").Value = 0.071
$ws.Range("J9").Value = 0.073
$ws.Range("K9").Value = 0.063
$ws.Range("L9").Value = 0.072

# Trial 9 - Oxytocin (ratio)
$ws.Range("B10").Value = "Oxytocin "
$ws.Range("D10").Value = "ratio"
$ws.Range("H10").Value = 0.996
$ws.Range("I10").Value = 0.024
$ws.Range("J10").Value = 0.025
$ws.Range("K10").Value = 0.995
$ws.Range("L10").Value = 0.025

# Trial 10 - Video
$ws.Range("B11").Value = "Video "
$ws.Range("D11").Value = "diff"
$ws.Range("H11").Value = -0.105
$ws.Range("I11").Value = 0.059
$ws.Range("J11").Value = 0.062
$ws.Range("K11").Value = -0.094
$ws.Range("L11").Value = 0.061

# Trial 11 - active
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "active "
$ws.Range("D12").Value = "diff"
$ws.Range("E12").Value = "diff"
$ws.Range("H12").Value = 0.004
$ws.Range("I12").Value = 0.187
$ws.Range("J12").Value = 0.395
$ws.Range("K12").Value = -0.214
$ws.Range("L12").Value = 0.253
$ws.Range("R12").Value = -0.027
$ws.Range("S12").Value = 0.231
$ws.Range("T12").Value = 0.426
$ws.Range("U12").Value = 0.615
$ws.Range("V12").Value = 0.371

# Trial 12 - no data yet
$ws.Range("A13").Value = 12

# Trial 13 - IV infusion
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "IV infusion"
$ws.Range("D14").Value = "diff"
$ws.Range("E14").Value = "diff"
$ws.Range("H14").Value = -42.364
$ws.Range("I14").Value = 29.871
$ws.Range("J14").Value = 30.12
$ws.Range("K14").Value = -42.83
$ws.Range("L14").Value = 30.534
$ws.Range("R14").Value = -0.03
$ws.Range("S14").Value = 0.037
$ws.Range("T14").Value = 0.038
$ws.Range("U14").Value = -0.033
$ws.Range("V14").Value = 0.038

# --- Placeholder rows for future trials (14 through 30), Trial_No only ---
for ($i = 15; $i -le 31; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# --- Restore the selection to where the author left off ---
$ws.Range("V17").Select()
